$d = $word.ActiveDocument

# --- Step 1: delete the "License Information" Heading2 paragraph entirely ---
$licInfoPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "License Information") {
        $licInfoPara = $para
        break
    }
}
if ($licInfoPara -ne $null) {
    $licInfoPara.Range.Delete()
}

# --- Step 2: locate the paragraph that begins the big license/attribution block ---
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Translation Questions (unfoldingWord) is based on")) {
        $targetPara = $para
        break
    }
}

$oldStart = $targetPara.Range.Start
$oldEnd = $targetPara.Range.End
$oldLen = $oldEnd - $oldStart

# Insert the replacement runs (with correct bold/plain formatting) right at the
# start of the paragraph, immediately after the existing leading empty run.
$insertionPoint = $d.Range($oldStart, $oldStart)
$insertionPoint.InsertXML('<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve">unfoldingWord® Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve"> © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. </w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve">unfoldingWord® Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve"> has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from </w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve">unfoldingWord® Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve"> © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Re-measure the paragraph after insertion, then delete all of the old
# (pre-existing) run content that now follows the newly inserted runs,
# leaving only the original trailing empty run in place.
$para2 = $d.Paragraphs.Item($targetPara.Index)
$newContentLen = $para2.Range.End - $para2.Range.Start - $oldLen
$deleteStart = $para2.Range.Start + $newContentLen
$deleteEnd = $deleteStart + $oldLen
$oldRange = $d.Range($deleteStart, $deleteEnd)
$oldRange.Delete()

# --- Step 3: delete the "This PDF version is provided under the same license." paragraph ---
$pdfPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "This PDF version is provided under the same license.") {
        $pdfPara = $para
        break
    }
}
if ($pdfPara -ne $null) {
    $pdfPara.Range.Delete()
}

Write-Host "Done."
